$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 4229.8887
$ws.Cells.Item(98, 9).Value = 4410.875
$ws.Cells.Item(98, 10).Value = 2782
$ws.Cells.Item(98, 11).Value = 4410.875
$ws.Cells.Item(98, 12).Value = 2782
$ws.Cells.Item(98, 13).Value = -2912.875
$ws.Cells.Item(98, 14).Value = -5778
$ws.Cells.Item(103, 8).Value = 614.25
$ws.Cells.Item(103, 10).Value = 595.8570999999999
$ws.Cells.Item(103, 12).Value = 1787.5713
$ws.Cells.Item(103, 14).Value = -2959.5713
$ws.Cells.Item(107, 8).Value = 2000.1
$ws.Cells.Item(107, 10).Value = 1976.8334
$ws.Cells.Item(107, 12).Value = 1976.8334
$ws.Cells.Item(107, 14).Value = -5816.8334
$ws.Cells.Item(113, 8).Value = 6384.1963
$ws.Cells.Item(113, 9).Value = 8571.286
$ws.Cells.Item(113, 10).Value = 6036.25
$ws.Cells.Item(113, 11).Value = 8571.286
$ws.Cells.Item(113, 12).Value = 6036.25
$ws.Cells.Item(113, 13).Value = -5317.286
$ws.Cells.Item(113, 14).Value = -12544.25
$ws.Cells.Item(122, 8).Value = 4229.8887
$ws.Cells.Item(122, 9).Value = 4410.875
$ws.Cells.Item(122, 10).Value = 2782
$ws.Cells.Item(122, 11).Value = 13232.625
$ws.Cells.Item(122, 12).Value = 8346
$ws.Cells.Item(122, 13).Value = -10782.625
$ws.Cells.Item(122, 14).Value = -13246
$ws.Cells.Item(132, 8).Value = 1163.5098
$ws.Cells.Item(132, 9).Value = 1182.0889
$ws.Cells.Item(132, 10).Value = 1024.1666
$ws.Cells.Item(132, 11).Value = 3546.2667
$ws.Cells.Item(132, 12).Value = 3072.4998
$ws.Cells.Item(132, 13).Value = -1016.2667
$ws.Cells.Item(132, 14).Value = -8132.4998
$ws.Cells.Item(135, 8).Value = 1028.7115
$ws.Cells.Item(135, 9).Value = 1072
$ws.Cells.Item(135, 11).Value = 9648
$ws.Cells.Item(135, 13).Value = -7113
$ws.Cells.Item(137, 8).Value = 1993.0968
$ws.Cells.Item(137, 9).Value = 2341.0527
$ws.Cells.Item(137, 10).Value = 1442.1666
$ws.Cells.Item(137, 11).Value = 7023.158100000001
$ws.Cells.Item(137, 12).Value = 4326.4998
$ws.Cells.Item(137, 13).Value = -4473.158100000001
$ws.Cells.Item(137, 14).Value = -9426.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 39.1
$ws.Cells.Item(5, 9).Value = 39.1
$ws.Cells.Item(5, 11).Value = 39.1
$ws.Cells.Item(5, 13).Value = 72.90000000000001
$ws.Cells.Item(32, 8).Value = 1751.0834
$ws.Cells.Item(32, 9).Value = 1052.9493
$ws.Cells.Item(32, 11).Value = 1052.9493
$ws.Cells.Item(32, 13).Value = -765.9493
$ws.Cells.Item(50, 8).Value = 2299.5715
$ws.Cells.Item(50, 9).Value = 1874.75
$ws.Cells.Item(50, 10).Value = 2866
$ws.Cells.Item(50, 11).Value = 1874.75
$ws.Cells.Item(50, 12).Value = 2866
$ws.Cells.Item(50, 13).Value = -1160.75
$ws.Cells.Item(50, 14).Value = -4294
$ws.Cells.Item(74, 8).Value = 2615.913
$ws.Cells.Item(74, 9).Value = 2196.2
$ws.Cells.Item(74, 10).Value = 3402.875
$ws.Cells.Item(74, 11).Value = 2196.2
$ws.Cells.Item(74, 12).Value = 3402.875
$ws.Cells.Item(74, 13).Value = -1322.2
$ws.Cells.Item(74, 14).Value = -5150.875
$ws.Cells.Item(77, 8).Value = 2615.913
$ws.Cells.Item(77, 9).Value = 2196.2
$ws.Cells.Item(77, 10).Value = 3402.875
$ws.Cells.Item(77, 11).Value = 10981
$ws.Cells.Item(77, 12).Value = 17014.375
$ws.Cells.Item(77, 13).Value = -6613
$ws.Cells.Item(77, 14).Value = -25750.375
$ws.Cells.Item(106, 8).Value = 28333.334
$ws.Cells.Item(106, 10).Value = 28333.334
$ws.Cells.Item(106, 12).Value = 28333.334
$ws.Cells.Item(106, 14).Value = -30857.334
$ws.Cells.Item(110, 8).Value = 988.61536
$ws.Cells.Item(110, 9).Value = 1011.7273
$ws.Cells.Item(110, 10).Value = 861.5
$ws.Cells.Item(110, 11).Value = 1011.7273
$ws.Cells.Item(110, 12).Value = 861.5
$ws.Cells.Item(110, 13).Value = 1033.2727
$ws.Cells.Item(110, 14).Value = -4951.5
$ws.Cells.Item(122, 8).Value = 4985.7144
$ws.Cells.Item(122, 9).Value = 4893.5713
$ws.Cells.Item(122, 10).Value = 5170
$ws.Cells.Item(122, 11).Value = 14680.7139
$ws.Cells.Item(122, 12).Value = 15510
$ws.Cells.Item(122, 13).Value = -12230.7139
$ws.Cells.Item(122, 14).Value = -20410

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 39.1
$ws.Cells.Item(4, 9).Value = 39.1
$ws.Cells.Item(4, 11).Value = 39.1
$ws.Cells.Item(4, 13).Value = 75.90000000000001
$ws.Cells.Item(5, 8).Value = 15239.8
$ws.Cells.Item(5, 9).Value = 15239.8
$ws.Cells.Item(5, 11).Value = 15239.8
$ws.Cells.Item(5, 13).Value = -15126.8
$ws.Cells.Item(46, 8).Value = 15181.818
$ws.Cells.Item(46, 9).Value = 15000
$ws.Cells.Item(46, 10).Value = 19000
$ws.Cells.Item(46, 11).Value = 15000
$ws.Cells.Item(46, 12).Value = 19000
$ws.Cells.Item(46, 13).Value = -14702
$ws.Cells.Item(46, 14).Value = -19596
$ws.Cells.Item(94, 8).Value = 1506.9535
$ws.Cells.Item(94, 9).Value = 953.76
$ws.Cells.Item(94, 10).Value = 2275.2778
$ws.Cells.Item(94, 11).Value = 953.76
$ws.Cells.Item(94, 12).Value = 2275.2778
$ws.Cells.Item(94, 13).Value = -502.76
$ws.Cells.Item(94, 14).Value = -3177.2778

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 6046.5713
$ws.Cells.Item(16, 9).Value = 6231
$ws.Cells.Item(16, 11).Value = 6231
$ws.Cells.Item(16, 13).Value = -5944
$ws.Cells.Item(31, 8).Value = 1688
$ws.Cells.Item(31, 9).Value = 972
$ws.Cells.Item(31, 10).Value = 2097.1428
$ws.Cells.Item(31, 11).Value = 972
$ws.Cells.Item(31, 12).Value = 2097.1428
$ws.Cells.Item(31, 13).Value = -677
$ws.Cells.Item(31, 14).Value = -2687.1428
$ws.Cells.Item(34, 8).Value = 1688
$ws.Cells.Item(34, 9).Value = 972
$ws.Cells.Item(34, 10).Value = 2097.1428
$ws.Cells.Item(34, 11).Value = 972
$ws.Cells.Item(34, 12).Value = 2097.1428
$ws.Cells.Item(34, 13).Value = -770
$ws.Cells.Item(34, 14).Value = -2501.1428
$ws.Cells.Item(58, 8).Value = 5107.3335
$ws.Cells.Item(58, 9).Value = 6104.8
$ws.Cells.Item(58, 10).Value = 120
$ws.Cells.Item(58, 11).Value = 6104.8
$ws.Cells.Item(58, 12).Value = 120
$ws.Cells.Item(58, 13).Value = -5901.8
$ws.Cells.Item(58, 14).Value = -526
$ws.Cells.Item(113, 8).Value = 6046.5713
$ws.Cells.Item(113, 9).Value = 6231
$ws.Cells.Item(113, 11).Value = 6231
$ws.Cells.Item(113, 13).Value = -4061
$ws.Cells.Item(132, 8).Value = 1689.72
$ws.Cells.Item(132, 9).Value = 1232.4736
$ws.Cells.Item(132, 10).Value = 3137.6667
$ws.Cells.Item(132, 11).Value = 3697.4208
$ws.Cells.Item(132, 12).Value = 9413.000100000001
$ws.Cells.Item(132, 13).Value = -1167.4208
$ws.Cells.Item(132, 14).Value = -14473.0001
$ws.Cells.Item(134, 8).Value = 1156.3158
$ws.Cells.Item(134, 9).Value = 1115
$ws.Cells.Item(134, 11).Value = 3345
$ws.Cells.Item(134, 13).Value = -810
$ws.Cells.Item(136, 8).Value = 5107.3335
$ws.Cells.Item(136, 9).Value = 6104.8
$ws.Cells.Item(136, 10).Value = 120
$ws.Cells.Item(136, 11).Value = 18314.4
$ws.Cells.Item(136, 12).Value = 360
$ws.Cells.Item(136, 13).Value = -15764.4
$ws.Cells.Item(136, 14).Value = -5460

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 109299.2
$ws.Cells.Item(37, 10).Value = 109299.2
$ws.Cells.Item(37, 12).Value = 327897.6
$ws.Cells.Item(37, 14).Value = -328121.6
$ws.Cells.Item(86, 14).Value = $null
$ws.Cells.Item(86, 8).Value = 99
$ws.Cells.Item(86, 9).Value = 99
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 297
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 889
$ws.Cells.Item(89, 14).Value = $null
$ws.Cells.Item(89, 8).Value = 99
$ws.Cells.Item(89, 9).Value = 99
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 891
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = 5037
$ws.Cells.Item(107, 8).Value = 375.875
$ws.Cells.Item(107, 9).Value = 351.66666
$ws.Cells.Item(107, 10).Value = 390.4
$ws.Cells.Item(107, 11).Value = 1054.99998
$ws.Cells.Item(107, 12).Value = 1171.2
$ws.Cells.Item(107, 13).Value = 865.0000199999999
$ws.Cells.Item(107, 14).Value = -5011.2
$ws.Cells.Item(132, 8).Value = 1590.3636
$ws.Cells.Item(132, 9).Value = 1549.4
$ws.Cells.Item(132, 11).Value = 13944.6
$ws.Cells.Item(132, 13).Value = -11414.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(54, 8).Value = 34987
$ws.Cells.Item(54, 10).Value = 34987
$ws.Cells.Item(54, 12).Value = 34987
$ws.Cells.Item(54, 14).Value = -35767
$ws.Cells.Item(107, 8).Value = 1457.25
$ws.Cells.Item(107, 9).Value = 2665.5
$ws.Cells.Item(107, 10).Value = 732.3
$ws.Cells.Item(107, 11).Value = 2665.5
$ws.Cells.Item(107, 12).Value = 732.3
$ws.Cells.Item(107, 13).Value = -745.5
$ws.Cells.Item(107, 14).Value = -4572.3
$ws.Cells.Item(113, 8).Value = 2400
$ws.Cells.Item(113, 10).Value = 2400
$ws.Cells.Item(113, 12).Value = 2400
$ws.Cells.Item(113, 14).Value = -6740
$ws.Cells.Item(132, 8).Value = 2042.5
$ws.Cells.Item(132, 9).Value = 2054.0476
$ws.Cells.Item(132, 10).Value = 1800
$ws.Cells.Item(132, 11).Value = 6162.1428
$ws.Cells.Item(132, 12).Value = 5400
$ws.Cells.Item(132, 13).Value = -3632.1428
$ws.Cells.Item(132, 14).Value = -10460

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6224.6875
$ws.Cells.Item(7, 9).Value = 3854.111
$ws.Cells.Item(7, 10).Value = 9272.571
$ws.Cells.Item(7, 11).Value = 3854.111
$ws.Cells.Item(7, 12).Value = 9272.571
$ws.Cells.Item(7, 13).Value = -3742.111
$ws.Cells.Item(7, 14).Value = -9496.571
$ws.Cells.Item(46, 8).Value = 1642.25
$ws.Cells.Item(46, 9).Value = 1106.3334
$ws.Cells.Item(46, 10).Value = 3250
$ws.Cells.Item(46, 11).Value = 1106.3334
$ws.Cells.Item(46, 12).Value = 3250
$ws.Cells.Item(46, 13).Value = -918.3334
$ws.Cells.Item(46, 14).Value = -3626
$ws.Cells.Item(55, 8).Value = 500.63635
$ws.Cells.Item(55, 9).Value = 247.8077
$ws.Cells.Item(55, 10).Value = 1439.7142
$ws.Cells.Item(55, 11).Value = 247.8077
$ws.Cells.Item(55, 12).Value = 1439.7142
$ws.Cells.Item(55, 13).Value = -74.80770000000001
$ws.Cells.Item(55, 14).Value = -1785.7142
$ws.Cells.Item(122, 8).Value = 6820.769
$ws.Cells.Item(122, 9).Value = 7000
$ws.Cells.Item(122, 10).Value = 6741.1113
$ws.Cells.Item(122, 11).Value = 21000
$ws.Cells.Item(122, 12).Value = 20223.3339
$ws.Cells.Item(122, 13).Value = -18550
$ws.Cells.Item(122, 14).Value = -25123.3339
$ws.Cells.Item(126, 8).Value = 6224.6875
$ws.Cells.Item(126, 9).Value = 3854.111
$ws.Cells.Item(126, 10).Value = 9272.571
$ws.Cells.Item(126, 11).Value = 11562.333
$ws.Cells.Item(126, 12).Value = 27817.713
$ws.Cells.Item(126, 13).Value = -9092.332999999999
$ws.Cells.Item(126, 14).Value = -32757.713
$ws.Cells.Item(132, 8).Value = 3244.0112
$ws.Cells.Item(132, 10).Value = 4905.8184
$ws.Cells.Item(132, 12).Value = 14717.4552
$ws.Cells.Item(132, 14).Value = -19777.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 195500
$ws.Cells.Item(17, 9).Value = 195500
$ws.Cells.Item(17, 11).Value = 195500
$ws.Cells.Item(17, 13).Value = -195328
$ws.Cells.Item(107, 8).Value = 815.3570999999999
$ws.Cells.Item(107, 9).Value = 543
$ws.Cells.Item(107, 10).Value = 1496.25
$ws.Cells.Item(107, 11).Value = 1629
$ws.Cells.Item(107, 12).Value = 4488.75
$ws.Cells.Item(107, 13).Value = 291
$ws.Cells.Item(107, 14).Value = -8328.75
$ws.Cells.Item(136, 8).Value = 6364.5293
$ws.Cells.Item(136, 9).Value = 6239.2856
$ws.Cells.Item(136, 11).Value = 18717.8568
$ws.Cells.Item(136, 13).Value = -16167.8568
